$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.303.24"
$ws.Range("E2").Value = "  +3.41%  "
$ws.Range("D3").Value = "2.068.91"
$ws.Range("E3").Value = "  +2.77%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.98"
$ws.Range("E5").Value = "  +1.92%  "
$ws.Range("E6").Value = "  +2.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.63"
$ws.Range("E7").Value = "  +11.54%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.389"
$ws.Range("E9").Value = "  +4.33%  "
$ws.Range("E10").Value = "  +4.63%  "
$ws.Range("E11").Value = "  +2.39%  "
$ws.Range("E12").Value = "  +6.81%  "
$ws.Range("E13").Value = "  +2.80%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.50"
$ws.Range("E14").Value = "  +8.55%  "
$ws.Range("E15").Value = "  +4.75%  "
$ws.Range("E16").Value = "  +3.68%  "
$ws.Range("D18").Value = "38.186.13"
$ws.Range("E18").Value = "  +3.30%  "
$ws.Range("E19").Value = "  +1.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.34"
$ws.Range("E20").Value = "  +2.16%  "
$ws.Range("E21").Value = "  +3.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.03"
$ws.Range("E22").Value = "  +2.45%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").Value = "  -1.37%  "
$ws.Range("E25").Value = "  +3.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.33"
$ws.Range("E26").Value = "  +4.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.71"
$ws.Range("E27").Value = "  +1.30%  "
$ws.Range("E28").Value = "  +5.77%  "
$ws.Range("E29").Value = "  +3.31%  "
$ws.Range("E30").Value = "  +2.28%  "
$ws.Range("E31").Value = "  +3.09%  "
$ws.Range("E32").Value = "  +4.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.66"
$ws.Range("E33").Value = "  +5.07%  "
$ws.Range("E34").Value = "  +9.95%  "
$ws.Range("E35").Value = "  +1.64%  "
$ws.Range("E36").Value = "  +0.61%  "
$ws.Range("E37").Value = "  +16.55%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.36"
$ws.Range("E38").Value = "  +6.84%  "
$ws.Range("E39").Value = "  +0.14%  "
$ws.Range("D40").Value = "1.527.43"
$ws.Range("E40").Value = "  +4.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.38"
$ws.Range("E41").Value = "  +9.27%  "
$ws.Range("E42").Value = "  +4.55%  "
$ws.Range("E43").Value = "  +3.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.88"
$ws.Range("E44").Value = "  +4.18%  "
$ws.Range("E45").Value = "  +2.51%  "
$ws.Range("E46").Value = "  +1.79%  "
$ws.Range("E47").Value = "  -1.92%  "
$ws.Range("E48").Value = "  +2.89%  "
$ws.Range("E49").Value = "  +2.10%  "
$ws.Range("E50").Value = "  +0.57%  "
$ws.Range("E51").Value = "  +2.97%  "
